$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.085.18"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").Value = "2.302.77"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'300.70"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").Value = "'97.99"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("D7").Value = "'0.509"
$ws.Range("E7").Value = "  +0.93%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D10").Value = "'33.76"
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("D11").Value = "'0.0795"
$ws.Range("E11").Value = "  +1.08%  "
$ws.Range("D12").Value = "'49.31"
$ws.Range("D14").Value = "'17.25"
$ws.Range("E14").Value = "  +13.60%  "
$ws.Range("E15").Value = "  +1.88%  "
$ws.Range("D16").Value = "2.660.05"
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("D17").Value = "2.302.01"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").Value = "'0.812"
$ws.Range("E18").Value = "  +3.49%  "
$ws.Range("D19").Value = "43.023.30"
$ws.Range("E19").Value = "  +1.65%  "
$ws.Range("D20").Value = "'11.73"
$ws.Range("E20").Value = "  +2.56%  "
$ws.Range("D21").Value = "0.0₃0903"
$ws.Range("E21").Value = "  +1.13%  "
$ws.Range("E22").Value = "  +1.43%  "
$ws.Range("D23").Value = "'67.83"
$ws.Range("E23").Value = "  +1.86%  "
$ws.Range("D24").Value = "'236.76"
$ws.Range("E24").Value = "  +0.93%  "
$ws.Range("D25").Value = "'2.06"
$ws.Range("E25").Value = "  +6.70%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("D28").Value = "'24.43"
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("E29").Value = "  +10.38%  "
$ws.Range("D30").Value = "'166.66"
$ws.Range("E30").Value = "  +1.71%  "
$ws.Range("D31").Value = "'34.14"
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("E35").Value = "  +6.73%  "
$ws.Range("D36").Value = "'2.39"
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("D37").Value = "'16.75"
$ws.Range("E37").Value = "  +4.02%  "
$ws.Range("E38").Value = "  +0.68%  "
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("E40").Value = "  +1.63%  "
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  -0.95%  "
$ws.Range("D44").Value = "1.981.27"
$ws.Range("E44").Value = "  +0.80%  "
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("D46").Value = "'9.94"
$ws.Range("E46").Value = "  +3.03%  "
$ws.Range("D47").Value = "'17.60"
$ws.Range("E47").Value = "  -1.17%  "
$ws.Range("E48").Value = "  +1.48%  "
$ws.Range("D49").Value = "2.528.45"
$ws.Range("E49").Value = "  +1.09%  "
$ws.Range("D50").Value = "'53.35"
$ws.Range("E50").Value = "  +1.26%  "
$ws.Range("D51").Value = "'4.60"
$ws.Range("E51").Value = "  -1.85%  "
